# Updated cryptos list on Mon Oct 23 09:52:05 UTC 2023 with GitHub Actions
#
# Refresh per-row Price (column D) and Volume(1h) (column E) figures for
# every coin, and apply the re-ranking of rows 13/14 (WrappedEther <->
# Polygon) and 16/17 (Polkadot <-> WrappedBTC), including their
# Coin/Link/Price/Volume cells.
#
# Some Price values are plain decimal numbers (e.g. "219.15"); Excel
# auto-detects those as numeric on assignment, which would round-trip
# through floating point and lose the exact text the source feed used
# (and the canonical file stores every Price/Volume cell as text). Force
# those cells to Text format first, write the literal string, then drop
# back to the default (unstyled) cell style so no stray formatting is
# introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.459.15"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").Value = "1.668.91"
$ws.Range("E3").Value = "  +2.23%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.525"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.35%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.51"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.263"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0631"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0903"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("D12").Value = "1.913.26"
$ws.Range("E12").Value = "  +2.43%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.609"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.23%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.655.40"
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "10.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.07%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "30.506.32"
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("D20").Value = "0.0₃0716"
$ws.Range("E20").Value = "  +2.24%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("E27").Value = "  +1.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0491"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("E31").Value = "  +2.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.17%  "
$ws.Range("D34").Value = "1.491.46"
$ws.Range("E34").Value = "  +4.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.74"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "83.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0177"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.588"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.97%  "
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.833"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0498"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.30%  "
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "50.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.06%  "
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "93.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.79%  "
$ws.Range("D51").Value = "0.0₆0109"
$ws.Range("E51").Value = "  -3.25%  "
